# "fixed template, fixed legal condition"
#
# The quote template's "Terms & Conditions" block used to show two extra
# rows with raw hyperlinked URLs (pointing at the General Terms of Service /
# Data Processing Agreement PDFs) right above the "Validity period" line.
# Those rows (and their hyperlinks) are removed, which naturally reflows
# everything below them up by 4 rows. The annual-fee formula in J30 is also
# fixed to use the commitment discount rate held in C8 instead of a
# hard-coded *1.2, and the print area / selection are refreshed to match
# the now-shorter sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the two hyperlinks that lived on A44 / A46 (the GENERAL TERMS
#        OF SERVICE / DATA PROCESSING AGREEMENT links) before the rows that
#        carry them are deleted, so no stale hyperlink survives pointing at
#        the wrong (shifted) cell. ---
$ws.Hyperlinks.Delete()

# --- 2. Remove the 4 rows (43:46) holding the two dated legal clause
#        headers + their raw URL text; this shifts rows 47:104 up to 43:100
#        and drops the now-unreferenced shared strings automatically. ---
$ws.Rows("43:46").Delete()

# --- 3. The G47:J48 block (old G51:J52) had picked up a redundant
#        alignment-only cell style; reset it back to the plain style used by
#        its neighbours so it matches the rest of the block again. ---
$ws.Range("F47").Copy()
$ws.Range("G47:I48").PasteSpecial(-4122)
$ws.Range("J53").Copy()
$ws.Range("J47:J48").PasteSpecial(-4122)
$ws.Range("A1").Select()

# --- 4. Fix the annual-price formula: it was hard-coded to *1.2 (a flat
#        20% uplift); use the commitment discount percentage in C8 instead. ---
$ws.Range("J30").Formula = '=J29*(1+$C$8/100)'

# --- 5. The printable area used to run to row 85; after deleting the 4
#        rows above, the equivalent end row is 81. ---
$ws.PageSetup.PrintArea = '$A$1:$J$81'

# --- 6. Refresh the view/selection state to where the author left off. ---
$ws.Range("J30").Select()
try {
    $excel.ActiveWindow.ScrollRow = 24
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # view-scroll position isn't critical; ignore if unsupported
}

$wb.Save()
